$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: rename the "Data" column header to "Date"
$ws.Range("D1").Value = "Date"

# Update the active selection cell (as captured by the saved workbook)
$ws.Range("I6").Select()
